$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Discount 25-26 Sum")

# --- Column E: add "Discount Type" header and "Percentage" labels for existing rows ---
$ws.Cells.Item(1,5).Value = "Discount Type"
$ws.Cells.Item(2,5).Value = "Percentage"
$ws.Cells.Item(3,5).Value = "Percentage"
$ws.Cells.Item(4,5).Value = "Percentage"
$ws.Cells.Item(5,5).Value = "Percentage"

# --- Row 6: new GEM-X3 discount entry (copy formatting from row 5 first) ---
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$ws.Cells.Item(6,1).Value = "GEM-X "
$ws.Cells.Item(6,2).Value = "Gem-X3 Discount"
$ws.Cells.Item(6,3).Value = 45
$ws.Cells.Item(6,4).Value = 46034
$ws.Cells.Item(6,5).Value = "Amount"

# --- Row 7: new Kbrand K-1 discount entry ---
$ws.Range("B5:C5").Copy()
$ws.Range("B7:C7").PasteSpecial(-4122)
$ws.Cells.Item(7,1).Value = "Kbrand"
$ws.Cells.Item(7,2).Value = "K-1 Discount"
$ws.Cells.Item(7,3).Value = 200
$ws.Cells.Item(7,4).Value = 46388
$ws.Cells.Item(7,5).Value = "Amount"

# --- Make the Discount sheet the active tab with B2 selected ---
$ws.Activate()
$ws.Range("B2").Select()
